$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7144740000000001
$ws.Range("N2").Value = 2.143422
$ws.Range("O2").Value = 0.138796410342318
$ws.Range("P2").Value = 0.138796410342318
$ws.Range("Q2").Value = 0.9404121130200001
$ws.Range("R2").Value = 8.463709017180001
$ws.Range("S2").Value = 0.138796410342318
$ws.Range("T2").Value = 0.138796410342318

# Row 3 updates
$ws.Range("O3").Value = 0.8044215857867821
$ws.Range("P3").Value = 0.8044215857867821
$ws.Range("Q3").Value = 5.45034126879
$ws.Range("R3").Value = 49.05307141911
$ws.Range("S3").Value = 0.8044215857867821
$ws.Range("T3").Value = 0.8044215857867821

# Row 4 updates
$ws.Range("M4").Value = 0.2847646666666667
$ws.Range("N4").Value = 0.8542940000000001
$ws.Range("O4").Value = 0.05531945672713084
$ws.Range("P4").Value = 0.05531945672713083
$ws.Range("Q4").Value = 0.3748157972066667
$ws.Range("R4").Value = 3.37334217486
$ws.Range("S4").Value = 0.05531945672713084
$ws.Range("T4").Value = 0.05531945672713083

# New Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf10"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.31623
$ws.Range("H5").Value = 3.94869
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.007528666666666667
$ws.Range("N5").Value = 0.022586
$ws.Range("O5").Value = 0.00146254714376898
$ws.Range("P5").Value = 0.00146254714376898
$ws.Range("Q5").Value = 0.009909456926666668
$ws.Range("R5").Value = 0.08918511234000001
$ws.Range("S5").Value = 0.00146254714376898
$ws.Range("T5").Value = 0.00146254714376898
